$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from the existing H1 header cell so I1/J1 match formatting
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Determine last used row in column H (data rows start at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, "H").End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
